$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared "password" value (column B, rows 2-1040) had a trailing
# asterisk ("SAS2023TH*") that needs to be removed, leaving "SAS2023TH".
# Setting the whole column range in one shot lets the engine fold it back
# into the shared-strings table the same way a user's find/replace would.
$ws.Range("B2:B1040").Value = "SAS2023TH"

# Move the sheet's active selection from C15 to C7.
$ws.Range("C7").Select()
